$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The daily price-history table is prepended with a new top row each day
# and every older row shifts down by one. Copy/paste (rather than
# re-typing through .Value) keeps the existing date cells' literal text
# storage intact instead of having Excel reinterpret strings like
# "2025-12-14" as date serials.
$ws.Range("A2:D25").Copy($ws.Range("A3:D26"))

# Build the new top row's date as literal text via a scratch cell, then
# paste just the value into A2 so we dodge Excel's automatic text->date
# conversion without leaving the scratch cell's formatting behind.
$scratch = $ws.Cells.Item(1, 10)
$scratch.NumberFormat = "@"
$scratch.Value = "2025-12-15"
$scratch.Copy()
$ws.Cells.Item(2, 1).PasteSpecial(-4163)
$scratch.Clear()

$ws.Cells.Item(2, 2).Value = 783.5
$ws.Cells.Item(2, 3).Value = 1112
$ws.Cells.Item(2, 4).Value = 3610
